# 自动更新Excel文件 - 2025-10-15 23:11:42
# For every data row, advance one day:
#   - decrement the "剩余" (remaining days, column E) counter by 1
#   - unless it has just hit 1, in which case the cycle restarts:
#       E is reset to the row's total day count (column D)
#       F ("开始时间", column F) is reset to the new cycle start date 20251016
# Rows whose start-date value isn't a clean 8-digit date (data-entry errors)
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$newCycleStart = 20251016

for ($r = 2; $r -le $lastRow; $r++) {

    $dCell = $ws.Cells.Item($r, 4)   # D: 总天 (total days)
    $eCell = $ws.Cells.Item($r, 5)   # E: 剩余 (days remaining)
    $fCell = $ws.Cells.Item($r, 6)   # F: 开始时间 (cycle start date)

    $d = $dCell.Value2
    $e = $eCell.Value2
    $f = $fCell.Value2

    if ($null -eq $d -or $null -eq $e -or $null -eq $f) {
        continue
    }

    $fStr = [string]$f
    if ($fStr.Length -ne 8) {
        # malformed/unparseable date (e.g. "202510929") - leave row untouched
        continue
    }

    if ($e -eq 1) {
        $eCell.Value = $d
        $fCell.Value = $newCycleStart
    } else {
        $eCell.Value = $e - 1
    }
}
